$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 5000
$ws.Range("F6").Value = 5000
$ws.Range("F8").Value = 138
$ws.Range("F9").Value = 190
$ws.Range("F12").Value = 159
$ws.Range("F13").Value = 8209
$ws.Range("F14").Value = 268
$ws.Range("F17").Value = 596
$ws.Range("F18").Value = 2497
$ws.Range("F19").Value = 6319
$ws.Range("F21").Value = 6
$ws.Range("F23").Value = 2512
$ws.Range("F25").Value = 11
$ws.Range("F26").Value = 6348
$ws.Range("F27").Value = 175
$ws.Range("F32").Value = 6790
$ws.Range("F33").Value = 5
$ws.Range("F35").Value = 220
$ws.Range("F39").Value = 23
$ws.Range("F40").Value = 38
$ws.Range("F41").Value = 42
$ws.Range("F46").Value = 58
$ws.Range("F47").Value = 497
$ws.Range("F48").Value = 2200
$ws.Range("F50").Value = 1114

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 152
$ws.Range("F5").Value = 47

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 5000
$ws.Range("F5").Value = 5000
$ws.Range("F7").Value = 138
$ws.Range("F8").Value = 190
$ws.Range("F11").Value = 159
$ws.Range("F12").Value = 8209
$ws.Range("F13").Value = 8209
$ws.Range("F14").Value = 268
$ws.Range("F16").Value = 596
$ws.Range("F17").Value = 2497
$ws.Range("F19").Value = 153
$ws.Range("F20").Value = 6319
$ws.Range("F22").Value = 47
$ws.Range("F23").Value = 6
$ws.Range("F24").Value = 2512
$ws.Range("F27").Value = 11
$ws.Range("F28").Value = 6348
$ws.Range("F29").Value = 175
$ws.Range("F34").Value = 6790
$ws.Range("F36").Value = 220
$ws.Range("F39").Value = 38
$ws.Range("F44").Value = 58
$ws.Range("F45").Value = 497
$ws.Range("F47").Value = 2200
$ws.Range("F50").Value = 1114
